$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.074275970458984
$ws.Range("C1").Value = 2.846532106399536
$ws.Range("D1").Value = 3.183358907699585
$ws.Range("E1").Value = 15
